# Move/resize the embedded chart ("Chart 1") and update the active
# selection on the sheet, per the target edit:
#   - Chart anchor "from": col 6/276225 off, row 1/38100 off
#                -> col 3/161926 off, row 1/38099 off
#   - Chart anchor "to":   col 18/104775 off, row 26/95250 off
#                -> col 22/419100 off, row 32/123824 off
#   - Selected cell: U11 -> V22

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reposition / resize the chart. Excel re-derives the drawing's
# twoCellAnchor (from/to col,colOff,row,rowOff) from the Left/Top/Width/
# Height (in points) of the ChartObject, so we target the new bounding
# box directly.
$co = $ws.ChartObjects().Item(1)
$co.Left = 200.86628967765748
$co.Top = 17.99992125984252
$co.Width = 1130.5624212598425
$co.Height = 471.75

# Update the selected cell/range shown in the sheet view.
$ws.Range("V22").Select()
